$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $find"
    }
    return $result
}

# --- Paragraph 1 (The Good, first paragraph) ---
Replace-Text " opinions and helped towards the completion of our project. Even though some of our meetings were intense and had some heated moments, each team member responded with poise and used it become a better member of the team. This is very important in a team setting. Overcoming adversity as a team and become a better unit after is a very important skill.  " " opinions and helped towards the completion of our project. Throughout the project, our group had some intense and heated moments. For example, we had a meeting where we argued about why some of the development work was necessary, such as the user account system. After having a team discussion, we all reached the same conclusion that this portion of the project added valuable functionality that went over-and-beyond the requirements of the project. Each team member responded to these meetings in a positive manner, using these moments to become a better member of the team. Overcoming adversity as a team and becoming a better unit after is an important skill for any team. Our team exemplified this."

# Reposition the _GoBack bookmark to sit inside "some" (after "had so")
$bmr = $d.Content
$bmr.Find.Execute("had so", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmr.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmr) | Out-Null

# --- Paragraph 2 (Another success...) ---
Replace-Text "Another success that we had as a team was our communication. One of the most important aspects of working in a team is the communication and we nailed it. Everyone was active in our online group chats and when we had meetings, everyone was in attendance. This made work on the project must less stressful and helped overcoming some of the issues we faced much easier. " "Another success that we had as a team was our communication. When working in a team, one of the most important aspects is the communication. There was never a moment where communication lacked or fell through. Each team member was active in our online group chats and everyone participated in meetings. This made work on the project less stress and helped overcome some of the issues we faced. "

# --- Paragraph 3 (During development...) ---
Replace-Text "During development, the team worked very hard to get work done. Everyone stepped in and contributed to get worked one when it mattered the most. Every teammate was able to use their strengths to benefit the team. This was very important become some of team was not familiar with the development stack we were using. This allowed those members of the team to take time and learn the development stack before contributing. " "During development, each team worked hard to finish the work they received. Everyone stepped up and contributed to get work done when we were in crunch time. Also, each teammate was able to use their strengths to benefit the team. For example, two of our team members were strong with node.js and another teammate was strong in AngularJS.  This allowed team members to take time and learn the technologies that they were not familiar with, as well as the ability to feel comfortable during the project."

# --- Paragraph (Angular 2 issues) ---
Replace-Text "Angular was not compatible with the Angular 2 version we were using. When we updated the Angular 2 version, we were then unable to load any webpages. After doing a few hours of investigation, one of the Angular 2 Javascript files had a 404 error, causing the application to fail when loading. " "Angular 2 was not compatible with the Angular 2 version we were using. When we updated the Angular 2 version, we were then unable to load any webpages. After doing a few hours of investigation, one of the Angular 2 Javascript files had a 404 error, causing the application to fail when loading. To alleviate this issue, we resorted to using Angular 1. Angular 1 allowed for the same functionality and was a tool that each teammate was familiar with. Although we were not able to accomplish or secondary goal of learning a new technology, Angular 1 was a perfect solution to Angular 2 failing to work. "
